# Apply "repull data, push all data, mean calculation" edits:
# update the dSF column (column F) values for specific rows to match
# the re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 1
    9  = 3
    19 = 4
    20 = 0
    27 = -2
    28 = 3
    30 = 0
    39 = 2
    40 = 0
    48 = 4
    49 = -1
    54 = 1
    56 = -4
    63 = -2
    66 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
